$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in order) mapped to their 1-based column index used below:
# C=3 D=4 E=5 F=6 G=7 H=8 I=9 L=12 M=13
$cols = @(3, 4, 5, 6, 7, 8, 9, 12, 13)

$newValues = @{
    2 = @(0.1382852495259073, 0.02903157483691032, 0.2080691495604583, 0.9606442962047055, 0.815313896557015, 0.8579786372760907, 0.8094353589091767, 0.4254645300077584, 10.94419732318642)
    3 = @(0.1342804142170735, 0.03086542026373706, 0.1895767653421387, 0.9805847396468295, 0.83770751258416, 0.8843039319907433, 0.8454207989169404, 0.3802221542334507, 9.626327495980604)
    4 = @(0.131908940825582, 0.03208420846612015, 0.1783528601004249, 0.9950466054534388, 0.8538623890222823, 0.9020316300492794, 0.8696134713979475, 0.3526737102124002, 8.814650732536165)
    5 = @(0.1309642598465359, 0.03260362793371385, 0.1738103245136458, 1.001485976137381, 0.8610358355447829, 0.9096424842088737, 0.8799887874304488, 0.3415017047533411, 8.483213416106594)
    6 = @(0.1308086998370186, 0.03269123472750834, 0.1730578882369045, 1.002587867407556, 0.8622622058597358, 0.910929424534217, 0.8817424893939076, 0.3396497713958695, 8.428136702285144)
    7 = @(0.1318961129721146, 0.03209112216626409, 0.1782914730711909, 0.9951312535711807, 0.8539567620013031, 0.902132715667932, 0.8697513189789206, 0.3525228259133257, 8.810183630120662)
    8 = @(0.136886031925485, 0.02964426853635871, 0.2016649345797106, 0.967052196611661, 0.8225271895037451, 0.8667270371882125, 0.8214012492559064, 0.409814735025634, 10.490286130198)
    9 = @(0.1473811489579191, 0.02560930081531509, 0.2486127407807572, 0.930129197361353, 0.7806436390220028, 0.8100041866765224, 0.7437221646446304, 0.5241837910725167, 13.76753160554597)
    10 = @(0.1555497643282138, 0.02314872059617912, 0.2839034006018437, 0.9148183178614602, 0.7628498897288836, 0.7764960231527169, 0.697789017588164, 0.6097366903998989, 16.16864197100853)
    11 = @(0.159371601920185, 0.02214760898804613, 0.3001605875285662, 0.9105882960652849, 0.7577830949615247, 0.7631191334730545, 0.6794660274276012, 0.6490598665051834, 17.26050881681698)
    12 = @(0.1608346186565939, 0.02178628982385789, 0.3063485582386392, 0.9093942430016284, 0.7563177539874175, 0.7583300370389736, 0.6729107683165907, 0.6640149974574001, 17.67399688551171)
    13 = @(0.1605188213066526, 0.02186330278898652, 0.3050144152514349, 0.9096330329103353, 0.756612886557491, 0.7593490283083781, 0.6743052975104646, 0.6607911809110476, 17.58494249934603)
    14 = @(0.1594916458626017, 0.02211752210837403, 0.3006690247962212, 0.9104818038800033, 0.7576533498445599, 0.7627195441054653, 0.6789189631005854, 0.6502889089901487, 17.29452580378234)
    15 = @(0.1588645418526085, 0.02227557803899671, 0.2980115566612227, 0.9110552408934893, 0.7583502454241398, 0.7648203275803667, 0.6817952824475171, 0.6438645262138891, 17.11664255048964)
    16 = @(0.1553021743767999, 0.02321659036062584, 0.2828452548671123, 0.9151510219795256, 0.7632434644014126, 0.7774084720583971, 0.6990393617443829, 0.6071754837686001, 16.09728439228888)
    17 = @(0.153144270044379, 0.02382469376149032, 0.2735950032937069, 0.9183743065932504, 0.7670333643023355, 0.7856145355933108, 0.7102862408677311, 0.5847758057009287, 15.47189133082321)
    18 = @(0.151913032394873, 0.02418553728659845, 0.2682934773146712, 0.920484791571468, 0.7694970856844918, 0.7905095980854497, 0.7169964339927475, 0.5719298072174581, 15.11213789560264)
    19 = @(0.1514978465368131, 0.0243095922471035, 0.2665016575506201, 0.9212429919947311, 0.7703794761480225, 0.7921968220304763, 0.7193094303931389, 0.5675866699463938, 14.99032168057232)
    20 = @(0.153372950959124, 0.02375880799040075, 0.2745777271451288, 0.9180045283847704, 0.7666004150311352, 0.7847228007446034, 0.7090639236840559, 0.5871563446999062, 15.53846946941815)
    21 = @(0.1597929196364589, 0.02204236278856087, 0.3019444894417376, 0.9102213150331124, 0.7573352908799507, 0.7617219736270329, 0.6775533073824107, 0.6533718850080277, 17.37982695635924)
    22 = @(0.1640809452385383, 0.0210245905847195, 0.3200164715410239, 0.9075183449149336, 0.7539311881524782, 0.7483048532913301, 0.6591997479989544, 0.6970254953904202, 18.5834301849431)
    23 = @(0.1617837166302252, 0.02155800932523277, 0.3103532021781348, 0.9087379497820365, 0.7554993502335492, 0.7553152876420199, 0.6687858541826301, 0.6736900300593049, 17.94100039499864)
    24 = @(0.1532695351726261, 0.02378856000647289, 0.2741333862762048, 0.918170904535387, 0.7667952654824148, 0.7851254024694754, 0.7096157736227084, 0.5860800041913308, 15.50837014259719)
    25 = @(0.1444635031345172, 0.02661526583543861, 0.2357812136348514, 0.9380987359491968, 0.7897645028533873, 0.8239461288816159, 0.7628316547192497, 0.4929986556314532, 12.88243172265288)
}

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
    }
}
